$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Hide the rows that are now collapsed in the log (rows 2-23) ---
$ws.Range("A2:A23").EntireRow.Hidden = $true

# --- Append two new work-log entries at the bottom ---
# Insert two blank rows directly above the current last row (row 35). This
# shifts the existing row 35 (date 45730, "Changed web app security
# scanning tool..." with the special "last row" styling) down to row 37,
# which is exactly where that styling needs to end up.
$ws.Rows.Item(35).Insert()
$ws.Rows.Item(35).Insert()

# The two freshly inserted rows (now rows 35 & 36) need the regular
# (non-emphasized) row styling used throughout the table, so copy the
# formatting from row 34.
$ws.Range("A34:C34").Copy($ws.Range("A35:C35"))
$ws.Range("A34:C34").Copy($ws.Range("A36:C36"))

# Row 35 keeps its original content (only its style changed above).
$ws.Range("A35").Value = 45730
$ws.Range("B35").Value = 4
$ws.Range("C35").Value = "Changed web app security scanning tool : nikto, tested Ollama tool for analysis report."

# Row 36 is a brand new entry.
$ws.Range("A36").Value = 45732
$ws.Range("B36").Value = 4
$ws.Range("C36").Value = "Generated AI integrated script, generated progress report"

# Row 37 is the new last entry; it already carries the emphasized "last
# row" style inherited from the original row 35, so only its values need
# to be set.
$ws.Range("A37").Value = 45734
$ws.Range("B37").Value = 4
$ws.Range("C37").Value = "Worked on the front-end part of my script"

# --- Update the view: scroll/selection moved to the new bottom of the table ---
$ws.Activate()
$ws.Range("C33").Select() | Out-Null
